$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 40, shifting the existing rows 40-52 down to 41-53.
$ws.Rows("40:40").Insert()

# Populate the newly inserted row 40 with the new data record (same constant
# columns as the surrounding rows, with its own date/price/origin values).
$ws.Range("A40").Value = 5
$ws.Range("B40").Value = "Macroferia Regional de Talca"
$ws.Range("C40").Value = "Maule"
$ws.Range("D40").Value = 44508
$ws.Range("E40").Value = 7
$ws.Range("F40").Value = 100112026
$ws.Range("G40").Value = "Haba"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 500
$ws.Range("K40").Value = 6000
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = 6000
$ws.Range("N40").Value = "$/saco 25 kilos"
$ws.Range("O40").Value = "Región del Maule"
$ws.Range("P40").Value = 240
$ws.Range("Q40").Value = 25
$ws.Range("R40").Value = "Hortaliza"
